# ---------------------------------------------------------------------------
# feat: add 2022-Q4 data
#
# 1. Insert a brand-new worksheet named "2022-Q4" right after "总计" (so it
#    becomes the second tab, pushing "2022-Q3", "2022-Q2", ... down by one),
#    and fill it with the quarter's fund-holding table.
# 2. Insert a new row at the top of the "总计" (summary) sheet's data for the
#    "2022-Q4" quarter (count=6, value=0.23), pushing every other quarter row
#    down by one.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1: "总计" summary sheet - insert the new 2022-Q4 row at row 2.
# ---------------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item(1)
$wsTotal.Rows(2).Insert()

# The inserted row inherits formatting from the row above/below; start clean
# so only the index cell (column A) carries the bold/boxed "index" style,
# matching the look of every other data row on this sheet.
$wsTotal.Range("A2:D2").ClearFormats()

$idxCell = $wsTotal.Cells.Item(2, 1)
$idxCell.Font.Bold = $true
$idxCell.Borders.LineStyle = 1
$idxCell.HorizontalAlignment = -4108
$idxCell.VerticalAlignment = -4160

$wsTotal.Cells.Item(2, 1).Value = 0
$wsTotal.Cells.Item(2, 2).Value = "2022-Q4"
$wsTotal.Cells.Item(2, 3).Value = 6
$wsTotal.Cells.Item(2, 4).Value = 0.23

# ---------------------------------------------------------------------------
# Step 2: brand-new "2022-Q4" worksheet, placed right before the old
# "2022-Q3" tab (i.e. right after "总计").
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($wb.Worksheets.Item(2))
$ws2.Name = "2022-Q4"

# Header row (B1:H1) - bold, boxed, centered.
$hdr = $ws2.Range("B1:H1")
$hdr.Font.Bold = $true
$hdr.Borders.LineStyle = 1
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160

$ws2.Cells.Item(1, "B").NumberFormat = "@"
$ws2.Cells.Item(1, "B").Value = "基金代码"
$ws2.Cells.Item(1, "C").NumberFormat = "@"
$ws2.Cells.Item(1, "C").Value = "基金名称"
$ws2.Cells.Item(1, "D").NumberFormat = "@"
$ws2.Cells.Item(1, "D").Value = "基金规模"
$ws2.Cells.Item(1, "E").NumberFormat = "@"
$ws2.Cells.Item(1, "E").Value = "股票总仓位"
$ws2.Cells.Item(1, "F").NumberFormat = "@"
$ws2.Cells.Item(1, "F").Value = "仓位占比"
$ws2.Cells.Item(1, "G").NumberFormat = "@"
$ws2.Cells.Item(1, "G").Value = "持有市值(亿元)"
$ws2.Cells.Item(1, "H").NumberFormat = "@"
$ws2.Cells.Item(1, "H").Value = "仓位排名"

# Data rows 2-7.
$ws2.Cells.Item(2, "A").Value = 0
$ws2.Cells.Item(2, "B").NumberFormat = "@"
$ws2.Cells.Item(2, "B").Value = "004317"
$ws2.Cells.Item(2, "C").NumberFormat = "@"
$ws2.Cells.Item(2, "C").Value = "前海开源沪港深裕鑫灵活配置混合C"
$ws2.Cells.Item(2, "D").NumberFormat = "@"
$ws2.Cells.Item(2, "D").Value = "2.88"
$ws2.Cells.Item(2, "E").NumberFormat = "@"
$ws2.Cells.Item(2, "E").Value = "90.85"
$ws2.Cells.Item(2, "F").NumberFormat = "@"
$ws2.Cells.Item(2, "F").Value = "3.04"
$ws2.Cells.Item(2, "G").NumberFormat = "@"
$ws2.Cells.Item(2, "G").Value = "0.0876"
$ws2.Cells.Item(2, "H").Value = 8

$ws2.Cells.Item(3, "A").Value = 1
$ws2.Cells.Item(3, "B").NumberFormat = "@"
$ws2.Cells.Item(3, "B").Value = "004316"
$ws2.Cells.Item(3, "C").NumberFormat = "@"
$ws2.Cells.Item(3, "C").Value = "前海开源沪港深裕鑫灵活配置混合A"
$ws2.Cells.Item(3, "D").NumberFormat = "@"
$ws2.Cells.Item(3, "D").Value = "2.30"
$ws2.Cells.Item(3, "E").NumberFormat = "@"
$ws2.Cells.Item(3, "E").Value = "90.85"
$ws2.Cells.Item(3, "F").NumberFormat = "@"
$ws2.Cells.Item(3, "F").Value = "3.04"
$ws2.Cells.Item(3, "G").NumberFormat = "@"
$ws2.Cells.Item(3, "G").Value = "0.0699"
$ws2.Cells.Item(3, "H").Value = 8

$ws2.Cells.Item(4, "A").Value = 2
$ws2.Cells.Item(4, "B").NumberFormat = "@"
$ws2.Cells.Item(4, "B").Value = "007107"
$ws2.Cells.Item(4, "C").NumberFormat = "@"
$ws2.Cells.Item(4, "C").Value = "太平 MSCI 香港价值增强指数A"
$ws2.Cells.Item(4, "D").NumberFormat = "@"
$ws2.Cells.Item(4, "D").Value = "1.02"
$ws2.Cells.Item(4, "E").NumberFormat = "@"
$ws2.Cells.Item(4, "E").Value = "92.16"
$ws2.Cells.Item(4, "F").NumberFormat = "@"
$ws2.Cells.Item(4, "F").Value = "5.30"
$ws2.Cells.Item(4, "G").NumberFormat = "@"
$ws2.Cells.Item(4, "G").Value = "0.0541"
$ws2.Cells.Item(4, "H").Value = 5

$ws2.Cells.Item(5, "A").Value = 3
$ws2.Cells.Item(5, "B").NumberFormat = "@"
$ws2.Cells.Item(5, "B").Value = "005255"
$ws2.Cells.Item(5, "C").NumberFormat = "@"
$ws2.Cells.Item(5, "C").Value = "浦银安盛港股通量化混合A"
$ws2.Cells.Item(5, "D").NumberFormat = "@"
$ws2.Cells.Item(5, "D").Value = "0.34"
$ws2.Cells.Item(5, "E").NumberFormat = "@"
$ws2.Cells.Item(5, "E").Value = "59.70"
$ws2.Cells.Item(5, "F").NumberFormat = "@"
$ws2.Cells.Item(5, "F").Value = "2.67"
$ws2.Cells.Item(5, "G").NumberFormat = "@"
$ws2.Cells.Item(5, "G").Value = "0.0091"
$ws2.Cells.Item(5, "H").Value = 8

$ws2.Cells.Item(6, "A").Value = 4
$ws2.Cells.Item(6, "B").NumberFormat = "@"
$ws2.Cells.Item(6, "B").Value = "013224"
$ws2.Cells.Item(6, "C").NumberFormat = "@"
$ws2.Cells.Item(6, "C").Value = "浦银安盛港股通量化混合C"
$ws2.Cells.Item(6, "D").NumberFormat = "@"
$ws2.Cells.Item(6, "D").Value = "0.17"
$ws2.Cells.Item(6, "E").NumberFormat = "@"
$ws2.Cells.Item(6, "E").Value = "59.70"
$ws2.Cells.Item(6, "F").NumberFormat = "@"
$ws2.Cells.Item(6, "F").Value = "2.67"
$ws2.Cells.Item(6, "G").NumberFormat = "@"
$ws2.Cells.Item(6, "G").Value = "0.0045"
$ws2.Cells.Item(6, "H").Value = 8

$ws2.Cells.Item(7, "A").Value = 5
$ws2.Cells.Item(7, "B").NumberFormat = "@"
$ws2.Cells.Item(7, "B").Value = "007108"
$ws2.Cells.Item(7, "C").NumberFormat = "@"
$ws2.Cells.Item(7, "C").Value = "太平 MSCI 香港价值增强指数C"
$ws2.Cells.Item(7, "D").NumberFormat = "@"
$ws2.Cells.Item(7, "D").Value = "0.00"
$ws2.Cells.Item(7, "E").NumberFormat = "@"
$ws2.Cells.Item(7, "E").Value = "92.16"
$ws2.Cells.Item(7, "F").NumberFormat = "@"
$ws2.Cells.Item(7, "F").Value = "5.30"
$ws2.Cells.Item(7, "G").Value = 0
$ws2.Cells.Item(7, "H").Value = 5

# Index column (A2:A7) - same bold/boxed/centered style as the header.
$idxCol = $ws2.Range("A2:A7")
$idxCol.Font.Bold = $true
$idxCol.Borders.LineStyle = 1
$idxCol.HorizontalAlignment = -4108
$idxCol.VerticalAlignment = -4160
